$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Standard Excel XlPasteType enum values
$xlPasteFormats = -4122
$xlPasteValues  = -4163

# --- Title text updates (Volume number, date range) ---
$ws.Range("A8").Value = "Volume 30   Number  26"
$ws.Range("C9").Value = "Report Covering the Week  6/26/2023  Through  7/2/2023"

# --- Cells that need a style/type change (numeric <-> text "no data" marker) ---
# G15 / D18: became the text marker "0" (same shared text already used by e.g. C14)
$ws.Range("C14").Copy()
$ws.Range("G15").PasteSpecial($xlPasteFormats)
$ws.Range("C14").Copy()
$ws.Range("G15").PasteSpecial($xlPasteValues)

$ws.Range("C14").Copy()
$ws.Range("D18").PasteSpecial($xlPasteFormats)
$ws.Range("C14").Copy()
$ws.Range("D18").PasteSpecial($xlPasteValues)

# H15 / E18: became the text marker "***.*" (same shared text already used by e.g. E14)
$ws.Range("E14").Copy()
$ws.Range("H15").PasteSpecial($xlPasteFormats)
$ws.Range("E14").Copy()
$ws.Range("H15").PasteSpecial($xlPasteValues)

$ws.Range("E14").Copy()
$ws.Range("E18").PasteSpecial($xlPasteFormats)
$ws.Range("E14").Copy()
$ws.Range("E18").PasteSpecial($xlPasteValues)

# C17 / D20: text marker -> ordinary numeric cell (borrow numeric formatting from I14)
$ws.Range("I14").Copy()
$ws.Range("C17").PasteSpecial($xlPasteFormats)
$ws.Range("C17").Value = 2

$ws.Range("I14").Copy()
$ws.Range("D20").PasteSpecial($xlPasteFormats)
$ws.Range("D20").Value = 5

# E20 / L28 / L29: text marker -> ordinary numeric % cell (borrow formatting from M14)
$ws.Range("M14").Copy()
$ws.Range("E20").PasteSpecial($xlPasteFormats)
$ws.Range("E20").Value = -20

$ws.Range("M14").Copy()
$ws.Range("L28").PasteSpecial($xlPasteFormats)
$ws.Range("L28").Value = 100

$ws.Range("M14").Copy()
$ws.Range("L29").PasteSpecial($xlPasteFormats)
$ws.Range("L29").Value = 100

# --- Simple same-style value updates ---
$ws.Range("N15").Value = -61.538461538461
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 9
$ws.Range("H16").Value = -11.111111111111
$ws.Range("I16").Value = 49
$ws.Range("J16").Value = 60
$ws.Range("K16").Value = -18.333333333333
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = 6.521739130434
$ws.Range("N16").Value = -83.557046979865
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 5
$ws.Range("G17").Value = 15
$ws.Range("H17").Value = -66.666666666666
$ws.Range("I17").Value = 55
$ws.Range("J17").Value = 54
$ws.Range("K17").Value = 1.851851851851
$ws.Range("L17").Value = 83.333333333333
$ws.Range("M17").Value = 52.777777777777
$ws.Range("N17").Value = 14.583333333333
$ws.Range("C18").Value = 1
$ws.Range("F18").Value = 4
$ws.Range("G18").Value = 18
$ws.Range("H18").Value = -77.777777777777
$ws.Range("I18").Value = 57
$ws.Range("K18").Value = -24
$ws.Range("L18").Value = 23.913043478260
$ws.Range("M18").Value = -1.724137931034
$ws.Range("N18").Value = -85.891089108910
$ws.Range("C19").Value = 16
$ws.Range("D19").Value = 20
$ws.Range("E19").Value = -20
$ws.Range("G19").Value = 85
$ws.Range("H19").Value = -32.941176470588
$ws.Range("I19").Value = 366
$ws.Range("J19").Value = 398
$ws.Range("K19").Value = -8.040201005025
$ws.Range("L19").Value = 18.831168831168
$ws.Range("M19").Value = 11.585365853658
$ws.Range("N19").Value = -62.804878048780
$ws.Range("C20").Value = 4
$ws.Range("F20").Value = 9
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = 50
$ws.Range("I20").Value = 48
$ws.Range("J20").Value = 28
$ws.Range("K20").Value = 71.428571428571
$ws.Range("L20").Value = 26.315789473684
$ws.Range("M20").Value = 182.352941176471
$ws.Range("N20").Value = -91.608391608391
$ws.Range("C21").Value = 25
$ws.Range("D21").Value = 29
$ws.Range("E21").Value = -13.793103448275
$ws.Range("F21").Value = 83
$ws.Range("G21").Value = 133
$ws.Range("H21").Value = -37.593984962406
$ws.Range("I21").Value = 581
$ws.Range("J21").Value = 622
$ws.Range("K21").Value = -6.591639871382
$ws.Range("L21").Value = 22.058823529411
$ws.Range("M21").Value = 18.571428571428
$ws.Range("N21").Value = -74.956896551724
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = -25
$ws.Range("I22").Value = 15
$ws.Range("J22").Value = 13
$ws.Range("K22").Value = 15.384615384615
$ws.Range("L22").Value = 36.363636363636
$ws.Range("M22").Value = 15.384615384615
$ws.Range("D23").Value = 2
$ws.Range("J23").Value = 13
$ws.Range("K23").Value = 69.230769230769
$ws.Range("L23").Value = 100
$ws.Range("C24").Value = 19
$ws.Range("D24").Value = 22
$ws.Range("E24").Value = -13.636363636363
$ws.Range("F24").Value = 90
$ws.Range("G24").Value = 94
$ws.Range("H24").Value = -4.255319148936
$ws.Range("I24").Value = 526
$ws.Range("J24").Value = 688
$ws.Range("K24").Value = -23.546511627907
$ws.Range("L24").Value = -25.495750708215
$ws.Range("M24").Value = 6.262626262626
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 3
$ws.Range("F25").Value = 25
$ws.Range("G25").Value = 23
$ws.Range("H25").Value = 8.695652173913
$ws.Range("I25").Value = 120
$ws.Range("J25").Value = 117
$ws.Range("K25").Value = 2.564102564102
$ws.Range("L25").Value = 46.341463414634
$ws.Range("M25").Value = -8.396946564885
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 100
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = -66.666666666666
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 9
$ws.Range("H27").Value = -44.444444444444
$ws.Range("I27").Value = 20
$ws.Range("J27").Value = 33
$ws.Range("K27").Value = -39.393939393939
$ws.Range("L27").Value = 0
$ws.Range("F30").Value = 1
$ws.Range("H30").Value = 0
